# Daily attendance processing - reorder "Recorded By" names in column G.
# For every data row, if the comma-separated list of recorders in column G
# ends with "System", reverse the order of the whole list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"
        if ($parts.Count -gt 1 -and $parts[$parts.Count - 1] -eq "System") {
            $reversedParts = $parts[($parts.Count - 1)..0]
            $newValue = $reversedParts -join ", "
            $cell.Value = $newValue
        }
    }
}
